# Insert a new title/author slide at the very front of the deck.
#
# Matches the target edit: a new slide (Section Header layout -> title +
# body placeholders) is inserted as the new first slide, pushing every
# existing slide down by one position. The new slide carries the paper
# title, the author list, and a small "name / student id" textbox.

function Emu([double]$v) {
    # Shape position/size properties on this COM surface are expressed in
    # points (same convention as real PowerPoint), so EMU values captured
    # from OOXML need to be converted (1 pt = 12700 EMU). A tiny nudge
    # compensates for the float32 storage under the property setter so the
    # value round-trips back to the exact same integer EMU amount.
    return ($v / 12700.0) + 0.00002
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Layout 3 == "Section Header" (title + body idx=1 placeholders), which is
# what the new slide's shapes use.
$layout = $master.CustomLayouts.Item(3)
$slide = $p.Slides.AddSlide(1, $layout)

# ---- Title placeholder ------------------------------------------------
$title = $slide.Shapes.Item(1)
$title.Left = Emu(676208)
$title.Top = Emu(1148691)
$title.Width = Emu(10515600)
$title.Height = Emu(2239692)
$title.TextFrame.AutoSize = 2

$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Artificial Neural Networks for Microwave Computer-Aided Design: The State of the Art"
$titleRange.Font.Bold = $true

# ---- Body / authors placeholder ---------------------------------------
$body = $slide.Shapes.Item(2)
$body.Left = Emu(676208)
$body.Top = Emu(3675063)
$body.Width = Emu(10515600)
$body.Height = Emu(1500187)
$body.TextFrame.AutoSize = 2

$bodyRange = $body.TextFrame.TextRange
$bodyRange.Text = "Authors:`rFeng Feng , Weicong Na , Jing Jin, Jianan Zhang , Wei Zhang, and Qi-Jun Zhang"
$bodyRange.Font.Bold = $true
$bodyRange.Font.Color.RGB = 0
$bodyRange.Font.NameFarEast = "+mj-ea"
$bodyRange.Font.NameComplexScript = "+mj-cs"

# Re-split the author line into the same run boundaries the source deck
# uses around the (non-English) proper nouns, keeping identical formatting
# on every run.
$authorsPara = $bodyRange.Paragraphs(2, 1)
$segments = @("Feng ", "Feng", " , ", "Weicong", " Na , Jing Jin, ", "Jianan", " Zhang , Wei Zhang, and Qi-Jun Zhang")
$pos = 1
foreach ($seg in $segments) {
    $run = $authorsPara.Characters($pos, $seg.Length)
    $run.Font.Bold = $true
    $run.Font.Color.RGB = 0
    $run.Font.NameFarEast = "+mj-ea"
    $run.Font.NameComplexScript = "+mj-cs"
    $pos += $seg.Length
}

# ---- Free-floating "name / student id" textbox ------------------------
$tb = $slide.Shapes.AddTextbox(1, (Emu 9280187), (Emu 5461930), (Emu 2200795), (Emu 830997))
$tb.Name = "TextBox 3"
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1

$tbRange = $tb.TextFrame.TextRange
$tbRange.Text = "Kolli SivaKrishna`r700765428"
$tbRange.Font.Size = 24

$tbPara1 = $tbRange.Paragraphs(1, 1)
$nameRun = $tbPara1.Characters(7, 11)
$nameRun.Font.Size = 24
